# Applies the coin-price/volume refresh described by the commit diff.
# For D-column values that look numeric (e.g. "300.61", "1.00"), a leading
# apostrophe is used so Excel stores them as text (quote-prefixed), matching
# the original inline-string ("text") cell type instead of silently coercing
# them to numbers (which would drop formatting like trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.034.47'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '2.298.98'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''300.61'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = '''99.53'
$ws.Range("E6").Value = '  +1.71%  '
$ws.Range("D7").Value = '''0.508'
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("D9").Value = '''0.510'
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("D10").Value = '''36.23'
$ws.Range("E10").Value = '  +6.98%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '''17.73'
$ws.Range("E13").Value = '  +2.24%  '
$ws.Range("D14").Value = '''6.88'
$ws.Range("E14").Value = '  +1.54%  '
$ws.Range("D15").Value = '2.656.70'
$ws.Range("E15").Value = '  -0.17%  '
$ws.Range("D16").Value = '2.305.27'
$ws.Range("E16").Value = '  -1.72%  '
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("D18").Value = '42.937.86'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = '''12.79'
$ws.Range("E19").Value = '  +8.87%  '
$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("E21").Value = '  +0.48%  '
$ws.Range("D22").Value = '''67.88'
$ws.Range("E22").Value = '  +0.14%  '
$ws.Range("D23").Value = '''235.54'
$ws.Range("E24").Value = '  +6.31%  '
$ws.Range("D26").Value = '''2.44'
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").Value = '''24.92'
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("D28").Value = '''169.74'
$ws.Range("E28").Value = '  +2.08%  '
$ws.Range("D29").Value = '''34.37'
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("E30").Value = '  -10.37%  '
$ws.Range("D31").Value = '''9.14'
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("D32").Value = '''1.00'
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("D33").Value = '''5.05'
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("D34").Value = '''17.71'
$ws.Range("E34").Value = '  +5.18%  '
$ws.Range("D35").Value = '''4.63'
$ws.Range("E35").Value = '  -0.93%  '
$ws.Range("E36").Value = '  +1.11%  '
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("E38").Value = '  +0.55%  '
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("E40").Value = '  +0.50%  '
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0291'
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.985.02'
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").Value = '''2.24'
$ws.Range("E44").Value = '  -5.44%  '
$ws.Range("E45").Value = '  +2.05%  '
$ws.Range("E46").Value = '  +1.07%  '
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").Value = '''55.51'
$ws.Range("E48").Value = '  +4.15%  '
$ws.Range("E49").Value = '  +3.68%  '
$ws.Range("D50").Value = '2.523.05'
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("D51").Value = '''70.74'
$ws.Range("E51").Value = '  +0.76%  '
